# Generate Report for Handoff
# Regenerates the localization-status report: the two tracked source files
# are swapped for a new pair of markdown files (fdc9c1bb.../ffff85af2807...),
# status moves from "Handed back: in sync with en-US" to "Ready for handoff",
# timestamps refresh, the zh-cn handoff xliff hash changes, and the "Latest
# Target File" / "Latest Handback File" / "Latest Handback DateTime" columns
# reset to their not-yet-handed-back defaults on the locale sheets.

$wb = $excel.ActiveWorkbook

$oldUuid1 = "1aa29009-39e0-4b33-a645-3f348e20e891"
$oldUuid2 = "603718cb-1111-4a69-ba0a-989b0d347a7d"
$newUuid1 = "fdc9c1bb-7730-40ed-81eb-6bcf496919e7"
$newUuid2 = "ffff85af2807-2b9e-4e92-ab25-8fd6df1e79ac"

$oldStatus = "Handed back: in sync with en-US"
$newStatus = "Ready for handoff"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "$newUuid1.md"
$ws1.Range("B2").Value = "e2e\$newUuid1.md"
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus
$ws1.Range("G2").Value = "2016-08-26 17:03:35"

$ws1.Range("A3").Value = "$newUuid2.md"
$ws1.Range("B3").Value = "e2e\$newUuid2.md"
$ws1.Range("E3").Value = $newStatus
$ws1.Range("F3").Value = $newStatus
$ws1.Range("G3").Value = "2016-08-26 17:03:35"

$ws1.Columns.Item(5).ColumnWidth = 16.38265482584637
$ws1.Columns.Item(6).ColumnWidth = 16.38265482584637

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/58e90c183d990db36743bb915e13c91c591a1d74/e2e/$oldUuid1.md", [Type]::Missing, [Type]::Missing, "e2e\$newUuid1.md")
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/58e90c183d990db36743bb915e13c91c591a1d74/e2e/$oldUuid2.md", [Type]::Missing, [Type]::Missing, "e2e\$newUuid2.md")

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "$newUuid1.md"
$ws2.Range("C2").Value = $newStatus
$ws2.Range("G2").Value = "$newUuid1.76e99ae8dbfe0cc7975c923e647324dfd4707f28.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-26 17:03:30"
$ws2.Range("I2").Value = ""
$ws2.Range("J2").Value = ""
$ws2.Range("K2").Value = "0001-01-01 00:00:00"

$ws2.Range("A3").Value = "$newUuid2.md"
$ws2.Range("C3").Value = $newStatus
$ws2.Range("F3").Value = "True"
$ws2.Range("G3").Value = "$newUuid1.76e99ae8dbfe0cc7975c923e647324dfd4707f28.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-26 17:03:30"
$ws2.Range("I3").Value = ""
$ws2.Range("J3").Value = ""
$ws2.Range("K3").Value = "0001-01-01 00:00:00"

$ws2.Range("I2").Style = "Normal"
$ws2.Range("I3").Style = "Normal"
$ws2.Range("J2").Style = "Normal"
$ws2.Range("J3").Style = "Normal"

$ws2.Columns.Item(3).ColumnWidth = 16.38265482584637
$ws2.Columns.Item(9).ColumnWidth = 17.817272004627068
$ws2.Columns.Item(10).ColumnWidth = 20.872143700009268

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/58e90c183d990db36743bb915e13c91c591a1d74/e2e/$oldUuid1.md", [Type]::Missing, [Type]::Missing, "$newUuid1.md")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f50982e97a8605bdae5c7efffa87fbdd16f8b0cd/e2e/$oldUuid1.md", [Type]::Missing, [Type]::Missing, "$newUuid2.md")

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "$newUuid1.md"
$ws3.Range("C2").Value = $newStatus
$ws3.Range("G2").Value = "$newUuid1.76e99ae8dbfe0cc7975c923e647324dfd4707f28.de-de.xlf"
$ws3.Range("H2").Value = "2016-08-26 17:03:35"
$ws3.Range("I2").Value = ""
$ws3.Range("J2").Value = ""
$ws3.Range("K2").Value = "0001-01-01 00:00:00"

$ws3.Range("A3").Value = "$newUuid2.md"
$ws3.Range("C3").Value = $newStatus
$ws3.Range("F3").Value = "True"
$ws3.Range("G3").Value = "$newUuid1.76e99ae8dbfe0cc7975c923e647324dfd4707f28.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-26 17:03:35"
$ws3.Range("I3").Value = ""
$ws3.Range("J3").Value = ""
$ws3.Range("K3").Value = "0001-01-01 00:00:00"

$ws3.Range("I2").Style = "Normal"
$ws3.Range("I3").Style = "Normal"
$ws3.Range("J2").Style = "Normal"
$ws3.Range("J3").Style = "Normal"

$ws3.Columns.Item(3).ColumnWidth = 16.38265482584637
$ws3.Columns.Item(9).ColumnWidth = 17.817272004627068
$ws3.Columns.Item(10).ColumnWidth = 20.872143700009268

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/58e90c183d990db36743bb915e13c91c591a1d74/e2e/$oldUuid1.md", [Type]::Missing, [Type]::Missing, "$newUuid1.md")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/be54bb146fe1e48c5d8f42b564d0df6bc16b5ee0/e2e/$oldUuid1.md", [Type]::Missing, [Type]::Missing, "$newUuid2.md")
